$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @{ D = "new price text (or $null if unchanged)"; E = "new volume text" }
$updates = @{
    2  = @{ D = "65.352.69"; E = "  +1.02%  " }
    3  = @{ D = "3.206.93";  E = "  -0.81%  " }
    4  = @{ D = $null;       E = "  -0.03%  " }
    5  = @{ D = "575.37";    E = "  -0.25%  " }
    6  = @{ D = "167.93";    E = "  -2.34%  " }
    7  = @{ D = "0.597";     E = "  -4.29%  " }
    8  = @{ D = $null;       E = "  +0.03%  " }
    9  = @{ D = "0.119";     E = "  -2.26%  " }
    10 = @{ D = "6.75";      E = "  -0.11%  " }
    11 = @{ D = $null;       E = "  +0.99%  " }
    12 = @{ D = "3.772.85";  E = "  -0.45%  " }
    13 = @{ D = $null;       E = "  -0.39%  " }
    14 = @{ D = "65.309.26"; E = "  +0.80%  " }
    15 = @{ D = "25.71";     E = "  -0.09%  " }
    16 = @{ D = "3.203.01";  E = "  -0.80%  " }
    17 = @{ D = $null;       E = "  -0.71%  " }
    18 = @{ D = "413.40";    E = "  -0.47%  " }
    19 = @{ D = "12.92";     E = "  +0.68%  " }
    20 = @{ D = "5.35";      E = "  -0.67%  " }
    21 = @{ D = "7.20";      E = "  -0.25%  " }
    22 = @{ D = $null;       E = "  +0.13%  " }
    23 = @{ D = "69.57";     E = "  -0.97%  " }
    24 = @{ D = $null;       E = "  -1.24%  " }
    25 = @{ D = "0.492";     E = "  -0.74%  " }
    26 = @{ D = $null;       E = "  -4.14%  " }
    27 = @{ D = "8.91";      E = "  -0.71%  " }
    28 = @{ D = $null;       E = "  -0.02%  " }
    29 = @{ D = $null;       E = "  -0.46%  " }
    30 = @{ D = "21.57";     E = "  -1.00%  " }
    31 = @{ D = "4.98";      E = "  -0.20%  " }
    32 = @{ D = "6.42";      E = "  +0.10%  " }
    33 = @{ D = $null;       E = "  -0.73%  " }
    34 = @{ D = "156.73";    E = "  -0.53%  " }
    35 = @{ D = "1.37";      E = "  -1.64%  " }
    36 = @{ D = "2.752.52";  E = "  -1.96%  " }
    37 = @{ D = $null;       E = "  -0.38%  " }
    38 = @{ D = "24.37";     E = "  -4.26%  " }
    39 = @{ D = "4.16";      E = "  -1.03%  " }
    40 = @{ D = "0.714";     E = "  -1.22%  " }
    41 = @{ D = "0.0634";    E = "  +0.84%  " }
    42 = @{ D = "5.71";      E = "  -0.97%  " }
    43 = @{ D = $null;       E = "  +0.44%  " }
    44 = @{ D = "297.94";    E = "  -1.32%  " }
    45 = @{ D = "21.69";     E = "  -1.40%  " }
    46 = @{ D = "0.0994";    E = "  -1.17%  " }
    47 = @{ D = "1.00";      E = "  +0.01%  " }
    48 = @{ D = "1.98";      E = "  -9.21%  " }
    49 = @{ D = "5.79";      E = "  -0.40%  " }
    50 = @{ D = "10.45";     E = "  +0.43%  " }
    51 = @{ D = "0.907";     E = "  -2.12%  " }
}

# Cell values such as "575.37" parse as a plain Excel number, which would
# silently change the cell from text to a numeric type. Force those through
# as text (temporary Text number-format, write, then restore the default
# "Normal" style) so the literal string content is preserved exactly as in
# the source data, matching the original inline-string text cells.
foreach ($row in ($updates.Keys | Sort-Object)) {
    $vals = $updates[$row]

    if ($null -ne $vals.D) {
        $dcell = $ws.Range("D$row")
        $dcell.NumberFormat = "@"
        $dcell.Value = $vals.D
        $dcell.Style = "Normal"
    }

    $ws.Range("E$row").Value = $vals.E
}
